# Update existing Cirbp-Trem1 LR-pair data with new TPM-based values
# and append new rows (row span grows from 11 to 16 rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Cirbp", "Trem1", "ECs", 3, 1, 8.052024666666666, 24.156074, 0.1791339256439199, 0.190543499292721, 1, 0.3333333333333333, 0.006995333333333333, 0.020986, 0.003071634989180631, 0.003071634989180632, 0.05632659655155555, 0.506939368964, 0.0005502340337571459, 0.000585280079388437),
  @("ECs", "Cirbp", "Trem1", "Inflammatory-Mac", 3, 1, 8.052024666666666, 24.156074, 0.1791339256439199, 0.190543499292721, 3, 1, 1.994843666666667, 5.984531, 0.8759313262859122, 0.8759313262859124, 16.06253041014378, 144.562773691294, 0.1569090170720807, 0.166903020050632),
  @("ECs", "Cirbp", "Trem1", "Resolving-Mac", 3, 1, 8.052024666666666, 24.156074, 0.1791339256439199, 0.190543499292721, 2, 0.6666666666666666, 0.2755583333333333, 0.826675, 0.120997038724907, 0.120997038724907, 2.218802497105556, 19.96922247395, 0.02167467453808199, 0.02305519916270066),
  @("FAPs", "Cirbp", "Trem1", "ECs", 3, 1, 11.67745566666667, 35.032367, 0.2597891290326613, 0.2763358729852725, 1, 0.3333333333333333, 0.006995333333333333, 0.020986, 0.003071634989180631, 0.003071634989180632, 0.08168769487355555, 0.7351892538620001, 0.0007979773785454841, 0.000848802936227338),
  @("FAPs", "Cirbp", "Trem1", "Inflammatory-Mac", 3, 1, 11.67745566666667, 35.032367, 0.2597891290326613, 0.2763358729852725, 3, 1, 1.994843666666667, 5.984531, 0.8759313262859122, 0.8759313262859124, 23.29469847943078, 209.652286314877, 0.227557436348241, 0.2420512477243652),
  @("FAPs", "Cirbp", "Trem1", "Resolving-Mac", 3, 1, 11.67745566666667, 35.032367, 0.2597891290326613, 0.2763358729852725, 2, 0.6666666666666666, 0.2755583333333333, 0.826675, 0.120997038724907, 0.120997038724907, 3.217820221080556, 28.960381989725, 0.03143371530587478, 0.03343582232468001),
  @("Inflammatory-Mac", "Cirbp", "Trem1", "ECs", 3, 1, 9.528830666666666, 28.586492, 0.2119885264612334, 0.2254907075621384, 1, 0.3333333333333333, 0.006995333333333333, 0.020986, 0.003071634989180631, 0.003071634989180632, 0.06665734679022221, 0.599916121112, 0.0006511513751831685, 0.0006926251470829621),
  @("Inflammatory-Mac", "Cirbp", "Trem1", "Inflammatory-Mac", 3, 1, 9.528830666666666, 28.586492, 0.2119885264612334, 0.2254907075621384, 3, 1, 1.994843666666667, 5.984531, 0.8759313262859122, 0.8759313262859124, 19.00852750613911, 171.076747555252, 0.1856873911405843, 0.1975143745400527),
  @("Inflammatory-Mac", "Cirbp", "Trem1", "Resolving-Mac", 3, 1, 9.528830666666666, 28.586492, 0.2119885264612334, 0.2254907075621384, 2, 0.6666666666666666, 0.2755583333333333, 0.826675, 0.120997038724907, 0.120997038724907, 2.625748697122222, 23.6317382741, 0.02564998394546582, 0.02728370787500275),
  @("MuSCs", "Cirbp", "Trem1", "ECs", 2, 1, 8.0746515, 16.149303, 0.1796373063025349, 0.1273859611772359, 1, 0.3333333333333333, 0.006995333333333333, 0.020986, 0.003071634989180631, 0.003071634989180632, 0.056484878793, 0.338909272758, 0.0005517802354010245, 0.0003912831754824034),
  @("MuSCs", "Cirbp", "Trem1", "Inflammatory-Mac", 2, 1, 8.0746515, 16.149303, 0.1796373063025349, 0.1273859611772359, 3, 1, 1.994843666666667, 5.984531, 0.8759313262859122, 0.8759313262859124, 16.1076674053155, 96.646004431893, 0.157349943960008, 0.111581353924182),
  @("MuSCs", "Cirbp", "Trem1", "Resolving-Mac", 2, 1, 8.0746515, 16.149303, 0.1796373063025349, 0.1273859611772359, 2, 0.6666666666666666, 0.2755583333333333, 0.826675, 0.120997038724907, 0.120997038724907, 2.2250375095875, 13.350225057525, 0.02173558210712579, 0.01541332407757152),
  @("Resolving-Mac", "Cirbp", "Trem1", "ECs", 3, 1, 7.616784666666667, 22.850354, 0.1694511125596505, 0.1802439589826321, 1, 0.3333333333333333, 0.006995333333333333, 0.020986, 0.003071634989180631, 0.003071634989180632, 0.05328194767155556, 0.479537529044, 0.0005204919662938081, 0.0005536436509994913),
  @("Resolving-Mac", "Cirbp", "Trem1", "Inflammatory-Mac", 3, 1, 7.616784666666667, 22.850354, 0.1694511125596505, 0.1802439589826321, 3, 1, 1.994843666666667, 5.984531, 0.8759313262859122, 0.8759313262859124, 15.19429465266378, 136.748651873974, 0.1484275377649981, 0.1578813300466805),
  @("Resolving-Mac", "Cirbp", "Trem1", "Resolving-Mac", 3, 1, 7.616784666666667, 22.850354, 0.1694511125596505, 0.1802439589826321, 2, 0.6666666666666666, 0.2755583333333333, 0.826675, 0.120997038724907, 0.120997038724907, 2.098868488105555, 18.88981639295, 0.02050308282835861, 0.02180898528495209)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value = $row[$j]
    }
}
